# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff" everywhere it
#   appears (Overview sheet's zh-cn/de-de status cells, and the Status column
#   on each per-locale sheet).
# - The "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
#   are refreshed to the new handoff run's time.
# - The Status column widens on every sheet so the longer "Ready for handoff"
#   text still fits.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

# --- Handoff timestamps ---
$overview.Range("G2").Value = "2016-08-29 12:42:19"
$dede.Range("H2").Value = "2016-08-29 12:42:19"
$zhcn.Range("H2").Value = "2016-08-29 12:42:15"

# --- Widen the Status columns to fit "Ready for handoff" ---
$overview.Columns.Item(5).ColumnWidth = 16.35
$overview.Columns.Item(6).ColumnWidth = 16.35
$zhcn.Columns.Item(3).ColumnWidth = 16.35
$dede.Columns.Item(3).ColumnWidth = 16.35
